$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.698.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.992.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.88%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.51%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.990.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.90%  "
$ws.Range("E9").Value = "  -3.26%  "
$ws.Range("E10").Value = "  -6.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.83%  "
$ws.Range("E13").Value = "  -7.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.03%  "
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.486.37"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.779.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.985.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.13%  "
$ws.Range("E19").Value = "  -7.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "433.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.31%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.02%  "
$ws.Range("E25").Value = "  -5.66%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("E28").Value = "  -7.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.54%  "
$ws.Range("E30").Value = "  -8.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -11.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0933"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.11%  "
$ws.Range("E34").Value = "  -5.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.953"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "49.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0666"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.51%  "
$ws.Range("E39").Value = "  -8.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.106"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "372.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.672.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.68%  "
$ws.Range("E44").Value = "  -8.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "120.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("E47").Value = "  -7.54%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "33.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.55%  "
$ws.Range("E50").Value = "  -4.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.21%  "
